$d = $word.ActiveDocument

# --- Split 1: "{m" -> "{" + "m" ---
# Locate the "{m" token that starts the "{m:commentblock ...}" field text.
$rng1 = $d.Content
$rng1.Find.Execute("{m:commentblock", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitStart1 = $rng1.Start + 1
$part1 = $d.Range($splitStart1, $splitStart1 + 1)
# Force Word to break the run at this character boundary by toggling a
# character-formatting property on the sub-range and then reverting it.
$part1.Bold = 1
$part1.Bold = 0

# --- Split 2: " some important comment}" -> " some important comment" + "}" ---
$rng2 = $d.Content
$rng2.Find.Execute(" some important comment}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitStart2 = $rng2.End - 1
$part2 = $d.Range($splitStart2, $rng2.End)
$part2.Bold = 1
$part2.Bold = 0
